# 02.10.2017 add headers tag & ovewview
#
# A1 held a short placeholder string; replace it with the real header /
# overview description text, and enable word-wrap on that cell so the
# long text is readable (adds a wrapText cell style, same as the target
# workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "The confirmation message that the system plays to the called party after the called party opts out to DNC."
$ws.Range("A1").WrapText = $true

# Give the sheet a wider base column (matches the sheetFormatPr tweak in
# the target file) — best-effort, harmless if the host doesn't persist it.
$ws.StandardWidth = 15
